# Refresh crypto price/volume snapshot cells to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''98.346.62'
$ws.Range("E2").Value = '  +2.00%  '

$ws.Range("D3").Value = '''3.302.25'
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").Value = '''255.88'
$ws.Range("E5").Value = '  +4.73%  '

$ws.Range("D6").Value = '''621.81'
$ws.Range("E6").Value = '  +1.40%  '

$ws.Range("D7").Value = '''1.45'
$ws.Range("E7").Value = '  +30.66%  '

$ws.Range("D8").Value = '''0.403'
$ws.Range("E8").Value = '  +6.33%  '

$ws.Range("D10").Value = '''0.912'
$ws.Range("E10").Value = '  +17.70%  '

$ws.Range("D11").Value = '''3.295.92'
$ws.Range("E11").Value = '  +0.58%  '

$ws.Range("D12").Value = '''0.199'
$ws.Range("E12").Value = '  +0.74%  '

$ws.Range("D13").Value = '''38.70'
$ws.Range("E13").Value = '  +11.83%  '

$ws.Range("D14").Value = '''97.491.16'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '''0.0000248'
$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").Value = '''3.928.14'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").Value = '''5.46'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '''3.298.57'
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").Value = '''15.13'
$ws.Range("E20").Value = '  +2.49%  '

$ws.Range("D21").Value = '''6.20'
$ws.Range("E21").Value = '  +8.03%  '

$ws.Range("D22").Value = '''481.87'
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("E23").Value = '  +3.19%  '

$ws.Range("D24").Value = '''0.0000203'
$ws.Range("E24").Value = '  +1.08%  '

$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("D26").Value = '''88.66'
$ws.Range("E26").Value = '  +1.46%  '

$ws.Range("D27").Value = '''11.82'
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("D28").Value = '''0.294'
$ws.Range("E28").Value = '  +24.07%  '

$ws.Range("D29").Value = '''3.481.54'
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").Value = '''0.187'
$ws.Range("E31").Value = '  +4.35%  '

$ws.Range("E32").Value = '  +9.84%  '

$ws.Range("D33").Value = '''9.92'
$ws.Range("E33").Value = '  +8.49%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '''27.66'
$ws.Range("E35").Value = '  +2.91%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '''7.17'
$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.147'
$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("D38").Value = '''1.94'
$ws.Range("E38").Value = '  +1.36%  '

$ws.Range("D39").Value = '''24.84'
$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("D40").Value = '''0.455'
$ws.Range("E40").Value = '  +2.47%  '

$ws.Range("D41").Value = '''488.18'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").Value = '''3.65'
$ws.Range("E42").Value = '  +6.35%  '

$ws.Range("E43").Value = '  -1.57%  '

$ws.Range("D44").Value = '''0.793'
$ws.Range("E44").Value = '  +2.54%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '''3.12'
$ws.Range("E46").Value = '  -2.25%  '

$ws.Range("D47").Value = '''158.15'
$ws.Range("E47").Value = '  -1.69%  '

$ws.Range("D48").Value = '''1.90'
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").Value = '''0.841'
$ws.Range("E49").Value = '  +8.50%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''7.17'
$ws.Range("E50").Value = '  +15.01%  '

$ws.Range("E51").Value = '  +3.50%  '
